# Apply crypto price/volume updates scraped on Wed Jul  3 15:28:18 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values (A1 ref -> target text). Ordered to match sheet layout.
$updates = [ordered]@{
    "D2" = "60.223.96"
    "E2" = "  -2.90%  "
    "D3" = "3.307.85"
    "E3" = "  -3.02%  "
    "D4" = "1.00"
    "E4" = "  +0.13%  "
    "D5" = "557.44"
    "E5" = "  -3.19%  "
    "D6" = "141.73"
    "E6" = "  -4.39%  "
    "E7" = "  +0.11%  "
    "D8" = "3.311.40"
    "E8" = "  -2.98%  "
    "D9" = "0.474"
    "E9" = "  -2.27%  "
    "D10" = "7.86"
    "E10" = "  -1.60%  "
    "E11" = "  -3.53%  "
    "D12" = "0.406"
    "E12" = "  -1.77%  "
    "D13" = "3.879.48"
    "E13" = "  -2.93%  "
    "E14" = "  +0.22%  "
    "D15" = "26.79"
    "E15" = "  -5.55%  "
    "D16" = "3.304.88"
    "E16" = "  -2.71%  "
    "D17" = "0.0000165"
    "E17" = "  -3.37%  "
    "D18" = "60.238.89"
    "E18" = "  -2.84%  "
    "D19" = "6.16"
    "E19" = "  -3.57%  "
    "D20" = "14.41"
    "E20" = "  -0.49%  "
    "D21" = "8.63"
    "E21" = "  -3.60%  "
    "D22" = "374.85"
    "E22" = "  -1.44%  "
    "D23" = "74.22"
    "E23" = "  -0.77%  "
    "D24" = "0.542"
    "E24" = "  -4.05%  "
    "D26" = "3.445.59"
    "E26" = "  -3.59%  "
    "E27" = "  -8.03%  "
    "E28" = "  -4.43%  "
    "D29" = "1.00"
    "E29" = "  +0.05%  "
    "D30" = "7.21"
    "E30" = "  -5.61%  "
    "D31" = "1.00"
    "E31" = "  -0.03%  "
    "D32" = "7.65"
    "E32" = "  -3.24%  "
    "D33" = "2.03"
    "E33" = "  -4.45%  "
    "D34" = "22.54"
    "E34" = "  -2.21%  "
    "D35" = "1.26"
    "E35" = "  -5.49%  "
    "D36" = "5.15"
    "E36" = "  -5.65%  "
    "B37" = "Monero"
    "C37" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D37" = "165.84"
    "E37" = "  -2.08%  "
    "E38" = "  -5.69%  "
    "B39" = "Aptos"
    "C39" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D39" = "6.74"
    "E39" = "  -2.31%  "
    "B40" = "RenzoRestakedETH"
    "C40" = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
    "D40" = "3.339.58"
    "E40" = "  -3.10%  "
    "B41" = "EnergySwap"
    "C41" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D41" = "26.76"
    "E41" = "  -11.85%  "
    "D42" = "0.0737"
    "D43" = "41.96"
    "E43" = "  -1.11%  "
    "D44" = "0.753"
    "E44" = "  -3.90%  "
    "D45" = "4.19"
    "E45" = "  -3.72%  "
    "D46" = "1.59"
    "E46" = "  -4.94%  "
    "E47" = "  -4.45%  "
    "D48" = "2.361.06"
    "E48" = "  -7.11%  "
    "E49" = "  -0.11%  "
    "D50" = "6.52"
    "E50" = "  -5.39%  "
    "D51" = "21.26"
    "E51" = "  -6.19%  "
}

# Column D sometimes holds values that look like plain numbers (e.g. "1.00",
# "557.44", "0.0000165"). The sheet stores Price/Volume as literal text, so
# force those specific cells to Text format before writing the new value --
# otherwise Excel would silently coerce them to numbers (dropping formatting
# like the trailing zero in "1.00") or convert tiny decimals to scientific
# notation.
$forceTextCells = @(
    "D4", "D5", "D6", "D9", "D10", "D12", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D50", "D51"
)

foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

